$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.398.24"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "1.776.85"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.74"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4243"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3613"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07139"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8378"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.43"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "1.794.33"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.445"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.236"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06875"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "78.86"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008681"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "26.403.93"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.072"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "2.002.64"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.07"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.789"
$ws.Range("E26").Value = "  -9.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.064"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.18"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.811"
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08845"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7234"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.737"
$ws.Range("E36").Value = "  -6.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.093"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05131"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01880"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1610"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4910"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.612"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.334"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.973"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.72"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.642"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06180"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4434"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.716"
$ws.Range("E51").Value = "  +1.99%  "
